$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RD Report")

# Row 9
$ws.Range("C9").Value = 25
$ws.Range("D9").Value = 625

# Row 10
$ws.Range("C10").Value = 77
$ws.Range("D10").Value = 1925

# Row 11
$ws.Range("C11").Value = 70
$ws.Range("D11").Value = 1400

# Row 12
$ws.Range("C12").Value = 18
$ws.Range("D12").Value = 540
$ws.Range("J12").Formula = "=300+200+100+75+50+35+25"

# Row 13
$ws.Range("J13").Value = 261

# Row 14
$ws.Range("J14").Formula = "=J12+J13"

# Row 25 - clear the blank placeholder string
$ws.Range("E25").Value = ""

# Row 31 - D31 value change
$ws.Range("D31").Value = 300

# Row 39 - D39 value change
$ws.Range("D39").Value = 1000

# Row 42 - D42 value change
$ws.Range("D42").Value = 1700

# sheet view changes
$ws.Range("A1").Select()
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("J15").Select()
